# Append the SSA raw/clean data row for August 26th to the historical log.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the first empty row right after the existing data (column A holds
# the "Fecha" / date-like text values).
$xlUp = -4162
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End($xlUp).Row
$newRow = $lastRow + 1

$dateCell = $ws.Cells.Item($newRow, 1)
# Force the date-looking string to be stored as plain text (matching the
# other "Fecha" entries, which live in the shared strings table) instead
# of letting Excel auto-convert it into a date serial number.
$dateCell.NumberFormat = "@"
$dateCell.Value = "2020-08-26"
$dateCell.Style = "Normal"

$ws.Cells.Item($newRow, 2).Value = 573888
$ws.Cells.Item($newRow, 3).Value = 635729
$ws.Cells.Item($newRow, 4).Value = 81466
$ws.Cells.Item($newRow, 5).Value = 62076
$ws.Cells.Item($newRow, 6).Value = 25.7
